$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.190.01'
$ws.Range('E2').Value = '  +1.07%  '

$ws.Range('D3').Value = '2.827.72'
$ws.Range('E3').Value = '  +3.08%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = '''354.31'
$ws.Range('E5').Value = '  +6.10%  '

$ws.Range('D6').Value = '''113.90'
$ws.Range('E6').Value = '  -2.26%  '

$ws.Range('E7').Value = '  +2.49%  '

$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  -0.06%  '

$ws.Range('D9').Value = '''0.606'
$ws.Range('E9').Value = '  +6.29%  '

$ws.Range('D10').Value = '''41.96'
$ws.Range('E10').Value = '  +1.33%  '

$ws.Range('D11').Value = '''0.0849'
$ws.Range('E11').Value = '  -0.92%  '

$ws.Range('E12').Value = '  +0.13%  '

$ws.Range('E13').Value = '  +1.26%  '

$ws.Range('D14').Value = '''7.79'
$ws.Range('E14').Value = '  +3.19%  '

$ws.Range('D15').Value = '3.257.47'
$ws.Range('E15').Value = '  +2.54%  '

$ws.Range('D16').Value = '2.822.99'
$ws.Range('E16').Value = '  +2.66%  '

$ws.Range('E17').Value = '  +1.67%  '

$ws.Range('D18').Value = '52.132.12'
$ws.Range('E18').Value = '  +1.13%  '

$ws.Range('E19').Value = '  +1.88%  '

$ws.Range('D20').Value = '''7.31'
$ws.Range('E20').Value = '  +7.14%  '

$ws.Range('D21').Value = '''13.83'
$ws.Range('E21').Value = '  +2.79%  '

$ws.Range('D22').Value = '0.0₃0997'
$ws.Range('E22').Value = '  +2.24%  '

$ws.Range('D23').Value = '''270.75'
$ws.Range('E23').Value = '  -2.76%  '

$ws.Range('D24').Value = '''69.61'
$ws.Range('E24').Value = '  +0.16%  '

$ws.Range('D25').Value = '''2.79'
$ws.Range('E25').Value = '  +5.05%  '

$ws.Range('D26').Value = '''26.71'
$ws.Range('E26').Value = '  -0.02%  '

$ws.Range('E27').Value = '  +0.06%  '

$ws.Range('E28').Value = '  +1.01%  '

$ws.Range('E29').Value = '  +1.36%  '

$ws.Range('E30').Value = '  +0.80%  '

$ws.Range('E31').Value = '  +1.58%  '

$ws.Range('D32').Value = '''33.88'
$ws.Range('E32').Value = '  -3.04%  '

$ws.Range('D33').Value = '''5.89'
$ws.Range('E33').Value = '  +6.12%  '

$ws.Range('D34').Value = '''0.0442'
$ws.Range('E34').Value = '  +28.22%  '

$ws.Range('E35').Value = '  +1.71%  '

$ws.Range('E36').Value = '  +0.12%  '

$ws.Range('E37').Value = '  +1.54%  '

$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').Value = '''18.47'
$ws.Range('E38').Value = '  -2.33%  '

$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '''4.88'
$ws.Range('E39').Value = '  -1.62%  '

$ws.Range('E40').Value = '  +1.16%  '

$ws.Range('E41').Value = '  +9.18%  '

$ws.Range('D42').Value = '''127.77'
$ws.Range('E42').Value = '  +0.50%  '

$ws.Range('D43').Value = '''23.42'
$ws.Range('E43').Value = '  +2.10%  '

$ws.Range('E44').Value = '  +1.99%  '

$ws.Range('D45').Value = '''2.29'
$ws.Range('E45').Value = '  +1.36%  '

$ws.Range('E46').Value = '  +1.55%  '

$ws.Range('D47').Value = '2.042.45'
$ws.Range('E47').Value = '  -2.29%  '

$ws.Range('E48').Value = '  +3.35%  '

$ws.Range('D49').Value = '''0.972'
$ws.Range('E49').Value = '  +12.29%  '

$ws.Range('D50').Value = '''5.70'
$ws.Range('E50').Value = '  +3.10%  '

$ws.Range('E51').Value = '  +1.17%  '

